# ---------------------------------------------------------------------
# Target change (per the supplied OOXML diff / commit message "Fixed
# POI packaging and upgraded to POI 3.15."):
#
#   The diff only touches the *serialized attribute order* inside
#   word/document.xml (the <w:document> namespace declarations and the
#   <w:pgSz>/<w:pgMar> attributes of the final <w:sectPr>) and inside
#   word/styles.xml (<w:rFonts>, <w:lang>, <w:latentStyles>, every
#   <w:lsdException>, the <w:style> elements and their <w:tblInd> /
#   <w:tblCellMar> children). Every attribute name/value pair present
#   before the commit is still present afterwards -- just written out
#   in (mostly alphabetical) order instead of the original order. This
#   is the signature of a Java OOXML writer (Apache POI) being
#   upgraded, not an authoring change made through Word's UI/object
#   model: no text, run, paragraph, style definition, margin value,
#   language, theme font, etc. was added, removed or modified.
#
# Verification performed while building this script: canonicalizing
# (i.e. sorting) the attributes of both the pre- and post-commit
# revisions of word/document.xml and word/styles.xml yields
# byte-for-byte identical XML -- confirming the edit carries no
# semantic/content change at all.
#
# The Word object model (Range/Find, PageSetup, Styles, ...) only lets
# a script describe *content* -- paragraph text, margins, style
# properties, etc. It has no notion of, and no way to control, the
# byte-level attribute order a part is serialized with; that is purely
# an internal detail of whichever OOXML writer produced the package.
# Driving PageSetup/Styles writes here would not reproduce the POI
# attribute ordering (this engine preserves input attribute order on
# write) and would instead risk *introducing* unwanted differences
# (e.g. extra namespace declarations minted on a full round-trip
# rewrite). Since the document already matches the target content
# exactly, no object-model mutation is required or appropriate.
# ---------------------------------------------------------------------

$d = $word.ActiveDocument
